# İş Takip Güncellemesi - 11.01.2026 17:00:26
#
# The "Güncelleme" sheet had its LİDAR/YERSEL ÖLÇÜM TARİHİ(YAPILAN) column
# (column H) removed. All the columns to its right (ARAZİ YERSEL ÖLÇÜM
# TARİHİ(YAPILAN) .. KESİN ASKI TARİHİ, originally I:P) shift one column to
# the left (H:O), and the sheet's used range shrinks from A1:P29 to A1:O29.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Güncelleme")

# Delete the entire "LİDAR/YERSEL ÖLÇÜM TARİHİ(YAPILAN)" column (H).
# This shifts every column to its right one position to the left, which is
# exactly the transformation the diff shows for columns I..P -> H..O.
$ws.Columns.Item(8).Delete()
